$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.404869
$ws.Range("H2").Value = 4.214607
$ws.Range("I2").Value = 0.6692718564235921
$ws.Range("J2").Value = 0.6692718564235923
$ws.Range("K2").Value = 3
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 1.404869
$ws.Range("N2").Value = 4.214607
$ws.Range("O2").Value = 0.6692718564235921
$ws.Range("P2").Value = 0.6692718564235923
$ws.Range("Q2").Value = 1.973656907161
$ws.Range("R2").Value = 17.762912164449
$ws.Range("S2").Value = 0.4479248178006813
$ws.Range("T2").Value = 0.4479248178006816

$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.404869
$ws.Range("H3").Value = 4.214607
$ws.Range("I3").Value = 0.6692718564235921
$ws.Range("J3").Value = 0.6692718564235923
$ws.Range("O3").Value = 0.3150411080808892
$ws.Range("P3").Value = 0.3150411080808893
$ws.Range("Q3").Value = 0.9290440843069999
$ws.Range("R3").Value = 8.361396758763
$ws.Range("S3").Value = 0.2108481472550423
$ws.Range("T3").Value = 0.2108481472550423

$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.404869
$ws.Range("H4").Value = 4.214607
$ws.Range("I4").Value = 0.6692718564235921
$ws.Range("J4").Value = 0.6692718564235923
$ws.Range("M4").Value = 0.03292866666666667
$ws.Range("N4").Value = 0.098786
$ws.Range("O4").Value = 0.01568703549551856
$ws.Range("P4").Value = 0.01568703549551856
$ws.Range("Q4").Value = 0.04626046301133333
$ws.Range("R4").Value = 0.416344167102
$ws.Range("S4").Value = 0.01049889136786849
$ws.Range("T4").Value = 0.01049889136786849

$ws.Range("I5").Value = 0.3150411080808892
$ws.Range("J5").Value = 0.3150411080808893
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 1.404869
$ws.Range("N5").Value = 4.214607
$ws.Range("O5").Value = 0.6692718564235921
$ws.Range("P5").Value = 0.6692718564235923
$ws.Range("Q5").Value = 0.9290440843069999
$ws.Range("R5").Value = 8.361396758763
$ws.Range("S5").Value = 0.2108481472550423
$ws.Range("T5").Value = 0.2108481472550423

$ws.Range("I6").Value = 0.3150411080808892
$ws.Range("J6").Value = 0.3150411080808893
$ws.Range("O6").Value = 0.3150411080808892
$ws.Range("P6").Value = 0.3150411080808893
$ws.Range("S6").Value = 0.09925089978083451
$ws.Range("T6").Value = 0.09925089978083454

$ws.Range("I7").Value = 0.3150411080808892
$ws.Range("J7").Value = 0.3150411080808893
$ws.Range("M7").Value = 0.03292866666666667
$ws.Range("N7").Value = 0.098786
$ws.Range("O7").Value = 0.01568703549551856
$ws.Range("P7").Value = 0.01568703549551856
$ws.Range("Q7").Value = 0.02177582605266667
$ws.Range("R7").Value = 0.195982434474
$ws.Range("S7").Value = 0.004942061045012407
$ws.Range("T7").Value = 0.004942061045012408

$ws.Range("G8").Value = 0.03292866666666667
$ws.Range("H8").Value = 0.098786
$ws.Range("I8").Value = 0.01568703549551856
$ws.Range("J8").Value = 0.01568703549551856
$ws.Range("K8").Value = 3
$ws.Range("L8").Value = 1
$ws.Range("M8").Value = 1.404869
$ws.Range("N8").Value = 4.214607
$ws.Range("O8").Value = 0.6692718564235921
$ws.Range("P8").Value = 0.6692718564235923
$ws.Range("Q8").Value = 0.04626046301133333
$ws.Range("R8").Value = 0.416344167102
$ws.Range("S8").Value = 0.01049889136786849
$ws.Range("T8").Value = 0.01049889136786849

$ws.Range("G9").Value = 0.03292866666666667
$ws.Range("H9").Value = 0.098786
$ws.Range("I9").Value = 0.01568703549551856
$ws.Range("J9").Value = 0.01568703549551856
$ws.Range("O9").Value = 0.3150411080808892
$ws.Range("P9").Value = 0.3150411080808893
$ws.Range("Q9").Value = 0.02177582605266667
$ws.Range("R9").Value = 0.195982434474
$ws.Range("S9").Value = 0.004942061045012407
$ws.Range("T9").Value = 0.004942061045012408

$ws.Range("G10").Value = 0.03292866666666667
$ws.Range("H10").Value = 0.098786
$ws.Range("I10").Value = 0.01568703549551856
$ws.Range("J10").Value = 0.01568703549551856
$ws.Range("M10").Value = 0.03292866666666667
$ws.Range("N10").Value = 0.098786
$ws.Range("O10").Value = 0.01568703549551856
$ws.Range("P10").Value = 0.01568703549551856
$ws.Range("Q10").Value = 0.001084297088444445
$ws.Range("R10").Value = 0.009758673795999999
$ws.Range("S10").Value = 0.0002460830826376592
$ws.Range("T10").Value = 0.0002460830826376592
